# fix(gui) step 1 and 2
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Step 1: bump the date in A1 by one day (45308 -> 45309)
$ws.Range("A1").Value = 45309

# Step 2: update the prices in column D for rows 33-38
$ws.Range("D33").Value = 3823.437
$ws.Range("D34").Value = 3823.437
$ws.Range("D35").Value = 5001.688
$ws.Range("D36").Value = 5001.725
$ws.Range("D37").Value = 5001.688
$ws.Range("D38").Value = 6304.761
